$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 0.20499999999999999
$ws.Range("J8").Value = 0.20699999999999999
$ws.Range("N8").Value = 0.20699999999999999
$ws.Range("O8").Value = 0.20699999999999999

$ws.Range("O12").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
